$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "25.481.90"
$ws.Cells.Item(2, 5).Value = "  +1.20%  "
$ws.Cells.Item(3, 4).Value = "1.661.48"
$ws.Cells.Item(3, 5).Value = "  +0.47%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9994"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "236.21"
$ws.Cells.Item(5, 5).Value = "  -1.45%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(6, 5).Value = "  +0.08%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4764"
$ws.Cells.Item(7, 5).Value = "  -1.08%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2611"
$ws.Cells.Item(8, 5).Value = "  -1.27%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06148"
$ws.Cells.Item(9, 5).Value = "  +2.11%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07072"
$ws.Cells.Item(10, 5).Value = "  -0.91%  "
$ws.Cells.Item(11, 4).Value = "1.657.97"
$ws.Cells.Item(11, 5).Value = "  +0.27%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "14.67"
$ws.Cells.Item(12, 5).Value = "  +0.61%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.5853"
$ws.Cells.Item(13, 5).Value = "  -6.44%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.347"
$ws.Cells.Item(14, 5).Value = "  -5.74%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "74.34"
$ws.Cells.Item(15, 5).Value = "  +0.91%  "
$ws.Cells.Item(16, 5).Value = "  +0.01%  "
$ws.Cells.Item(17, 5).Value = "  +0.08%  "
$ws.Cells.Item(18, 4).Value = "25.480.03"
$ws.Cells.Item(18, 5).Value = "  +1.28%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000006772"
$ws.Cells.Item(19, 5).Value = "  +2.73%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.36"
$ws.Cells.Item(20, 5).Value = "  -0.50%  "
$ws.Cells.Item(21, 4).Value = "1.868.29"
$ws.Cells.Item(21, 5).Value = "  +0.13%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.414"
$ws.Cells.Item(22, 5).Value = "  -1.61%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.607"
$ws.Cells.Item(23, 5).Value = "  -0.10%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.267"
$ws.Cells.Item(24, 5).Value = "  -0.88%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "132.60"
$ws.Cells.Item(25, 5).Value = "  +0.09%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "15.02"
$ws.Cells.Item(26, 5).Value = "  +1.39%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.392"
$ws.Cells.Item(27, 5).Value = "  -0.16%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "104.54"
$ws.Cells.Item(28, 5).Value = "  +2.01%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.676"
$ws.Cells.Item(29, 5).Value = "  -0.19%  "
$ws.Cells.Item(30, 5).Value = "  +4.32%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.639"
$ws.Cells.Item(31, 5).Value = "  +0.42%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.07635"
$ws.Cells.Item(32, 5).Value = "  -3.87%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.9993"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.04246"
$ws.Cells.Item(34, 5).Value = "  -7.43%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.614"
$ws.Cells.Item(35, 5).Value = "  -0.51%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.6078"
$ws.Cells.Item(36, 5).Value = "  +4.27%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.9466"
$ws.Cells.Item(37, 5).Value = "  -0.36%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.606"
$ws.Cells.Item(38, 5).Value = "  -1.51%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.8508"
$ws.Cells.Item(39, 5).Value = "  -0.05%  "
$ws.Cells.Item(40, 5).Value = "  +0.03%  "
$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.853"
$ws.Cells.Item(41, 5).Value = "  +0.87%  "
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.01476"
$ws.Cells.Item(42, 5).Value = "  -4.88%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "96.88"
$ws.Cells.Item(43, 5).Value = "  -2.55%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.3753"
$ws.Cells.Item(44, 5).Value = "  +0.62%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "4.689"
$ws.Cells.Item(45, 5).Value = "  -2.71%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.1114"
$ws.Cells.Item(46, 5).Value = "  -2.16%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "6.181"
$ws.Cells.Item(47, 5).Value = "  +1.27%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.05252"
$ws.Cells.Item(48, 5).Value = "  +1.38%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "29.43"
$ws.Cells.Item(49, 5).Value = "  -1.07%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.002"
$ws.Cells.Item(50, 5).Value = "  +0.09%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.9995"
$ws.Cells.Item(51, 5).Value = "  +0.09%  "
